$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CAP")
$ws.Range("E2").Value = 0.05529194511868137
$ws.Range("F2").Value = 0.1596163206428224
$ws.Range("G2").Value = 0.3564546721448534
$ws.Range("H2").Value = 0.727847612935295
$ws.Range("I2").Value = 1.428588679990762
$ws.Range("J2").Value = 2.750741012808266
$ws.Range("K2").Value = 5.245366878708023
$ws.Range("L2").Value = 17.04825401325482
$ws.Range("M2").Value = 15.71994196771707
$ws.Range("N2").Value = 13.57890037051744
$ws.Range("O2").Value = 8.1301892460665
$ws.Range("P2").Value = 4.867844624618641
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("E5").Value = 0.1437590573085715
$ws.Range("F5").Value = 0.4150024336713384
$ws.Range("G5").Value = 0.926782147576619
$ws.Range("H5").Value = 1.892403793631766
$ws.Range("I5").Value = 3.714330567975981
$ws.Range("J5").Value = 7.151926633301487
$ws.Range("K5").Value = 13.63795388464085
$ws.Range("L5").Value = 44.32546043446251
$ws.Range("M5").Value = 41.389714838094
$ws.Range("N5").Value = 35.42955872964427
$ws.Range("O5").Value = 21.21298555235359
$ws.Range("P5").Value = 20.58286074629813
$ws.Range("E6").Value = 0.2488137530340661
$ws.Range("F6").Value = 0.7182734428927011
$ws.Range("G6").Value = 1.60404602465184
$ws.Range("H6").Value = 3.275314258208827
$ws.Range("I6").Value = 6.428649059958429
$ws.Range("J6").Value = 12.37833455763719
$ws.Range("K6").Value = 23.6041509541861
$ws.Range("L6").Value = 76.71714305964667
$ws.Range("M6").Value = 158.1062696695804
$ws.Range("N6").Value = 203.2149114614893
$ws.Range("O6").Value = 149.5705981339387
$ws.Range("P6").Value = 117.3975917819989
$ws.Range("F7").Value = 0.1879926134035166
$ws.Range("G7").Value = 0.5426954901855964
$ws.Range("H7").Value = 1.211945885292502
$ws.Range("I7").Value = 2.474681883980002
$ws.Range("J7").Value = 4.857201511968591
$ws.Range("K7").Value = 9.352519443548099
$ws.Range("L7").Value = 30.63053261173651
$ws.Range("M7").Value = 18.33963134319247
$ws.Range("N7").Value = 10.9806147371833
$ws.Range("O7").Value = 6.574499658696964
$ws.Range("P7").Value = 3.936395802671987
$ws.Range("E8").Value = 0.8072623987327479
$ws.Range("F8").Value = 2.330398281385208
$ws.Range("G8").Value = 5.204238213314861
$ws.Range("H8").Value = 10.62657514885531
$ws.Range("I8").Value = 20.85739472786513
$ws.Range("J8").Value = 40.16081878700066
$ws.Range("K8").Value = 76.58235642913714
$ws.Range("L8").Value = 248.9045085935203
$ws.Range("M8").Value = 240.8263962367501
$ws.Range("N8").Value = 226.5323386248103
$ws.Range("O8").Value = 137.6293227581042
$ws.Range("P8").Value = 133.1592310797446
$ws.Range("E9").Value = 0.1941074700677905
$ws.Range("F9").Value = 0.5763497281931176
$ws.Range("G9").Value = 1.297561195654932
$ws.Range("H9").Value = 2.658336907240128
$ws.Range("I9").Value = 5.22583703831401
$ws.Range("J9").Value = 10.07017462234186
$ws.Range("K9").Value = 19.21042990178723
$ws.Range("L9").Value = 62.45673652383566
$ws.Range("M9").Value = 60.61942377799791
$ws.Range("N9").Value = 53.57984118059964
$ws.Range("O9").Value = 33.90046200218654
$ws.Range("P9").Value = 33.73680997185478
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.2488137530340661
$ws.Range("I10").Value = 0.7182734428927011
$ws.Range("J10").Value = 1.60404602465184
$ws.Range("K10").Value = 3.275314258208827
$ws.Range("L10").Value = 8.162511276248772
$ws.Range("M10").Value = 4.887197018039945
$ws.Range("N10").Value = 2.926145384036168
$ws.Range("O10").Value = 1.751991331004326
$ws.Range("P10").Value = 1.048981927097704
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.7796164261734072
$ws.Range("J11").Value = 2.250590121063797
$ws.Range("K11").Value = 2.250590121063797
$ws.Range("L11").Value = 1.34751144056587
$ws.Range("M11").Value = 0.8068048755131079
$ws.Range("N11").Value = 0.4830638817273195
$ws.Range("O11").Value = 0.2892281900020255
$ws.Range("P11").Value = 0.1731716012232691
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("E13").Value = 0.6026822017936269
$ws.Range("F13").Value = 1.739817895006764
$ws.Range("G13").Value = 3.885355926378902
$ws.Range("H13").Value = 7.933538980994713
$ws.Range("I13").Value = 15.5716166118993
$ws.Range("J13").Value = 29.98307703961009
$ws.Range("K13").Value = 57.17449897791743
$ws.Range("L13").Value = 70.06569857895894
$ws.Range("M13").Value = 41.95092191276472
$ws.Range("N13").Value = 25.11756658427699
$ws.Range("O13").Value = 15.03881493778619
$ws.Range("P13").Value = 9.004294025622512
$ws.Range("E14").Value = 2.051916826055484
$ws.Range("F14").Value = 6.127450715195469
$ws.Range("G14").Value = 13.81713366990761
$ws.Range("H14").Value = 28.57477634019261
$ws.Range("I14").Value = 57.19898843904972
$ws.Range("J14").Value = 111.2069103103838
$ws.Range("K14").Value = 210.3331808491975
$ws.Range("L14").Value = 559.6583565322301
$ws.Range("M14").Value = 582.6463016396496
$ws.Range("N14").Value = 571.8429409542847
$ws.Range("O14").Value = 374.098091810139
$ws.Range("P14").Value = 323.9071815611306

$ws = $wb.Worksheets.Item("CAP_NEW")
$ws.Range("E2").Value = 0.01105838902373627
$ws.Range("F2").Value = 0.02086487510482822
$ws.Range("G2").Value = 0.03936767030040619
$ws.Range("H2").Value = 0.07427858815808831
$ws.Range("I2").Value = 0.1512066024348298
$ws.Range("J2").Value = 0.2852953416683288
$ws.Range("K2").Value = 0.5382928434803577
$ws.Range("L2").Value = 1.255229658142431
$ws.Range("M2").Value = 0.04761811688909753
$ws.Range("N2").Value = 0.07105245497644007
$ws.Range("O2").Value = 0.7419664696302101
$ws.Range("P2").Value = 0.04254172941796646
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("E5").Value = 0.02875181146171431
$ws.Range("F5").Value = 0.05424867527255337
$ws.Range("G5").Value = 0.1023559427810561
$ws.Range("H5").Value = 0.1931243292110295
$ws.Range("I5").Value = 0.3931371663305572
$ws.Range("J5").Value = 0.7417678883376546
$ws.Range("K5").Value = 1.39956139304893
$ws.Range("L5").Value = 3.26359711117032
$ws.Range("M5").Value = 0.1052423468070983
$ws.Range("N5").Value = 0.1570353385780007
$ws.Range("O5").Value = 1.964263216657359
$ws.Range("P5").Value = 0.09402285797245463
$ws.Range("E6").Value = 0.04976275060681323
$ws.Range("F6").Value = 0.09389193797172701
$ws.Range("G6").Value = 0.1771545163518279
$ws.Range("H6").Value = 0.3342536467113974
$ws.Range("I6").Value = 0.6804297109567335
$ws.Range("J6").Value = 1.28382903750748
$ws.Range("K6").Value = 2.422317795661608
$ws.Range("L6").Value = 5.633072717490311
$ws.Range("M6").Value = 8.58706054710906
$ws.Range("N6").Value = 5.53603348448003
$ws.Range("O6").Value = 1.076878691517947
$ws.Range("P6").Value = 5.126847002201915
$ws.Range("F7").Value = 0.03759852268070332
$ws.Range("G7").Value = 0.07094057535641596
$ws.Range("H7").Value = 0.1338500790213811
$ws.Range("I7").Value = 0.2525471997375001
$ws.Range("J7").Value = 0.514102448278421
$ws.Range("K7").Value = 0.9700041616723176
$ws.Range("L7").Value = 2.190938116753216
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 1.038271059729859
$ws.Range("O7").Value = 0.03579872946725225
$ws.Range("P7").Value = 0.3722055585573534
$ws.Range("E8").Value = 0.1614524797465496
$ws.Range("F8").Value = 0.3046271765304921
$ws.Range("G8").Value = 0.5747679863859306
$ws.Range("H8").Value = 1.084467387108089
$ws.Range("I8").Value = 2.207616395548515
$ws.Range("J8").Value = 4.165311988357599
$ws.Range("K8").Value = 7.859075514813224
$ws.Range("L8").Value = 17.7431737022993
$ws.Range("M8").Value = 0.2474047986986314
$ws.Range("N8").Value = 3.528702400943838
$ws.Range("O8").Value = 9.98682507616795
$ws.Range("P8").Value = 0.2978358107747358
$ws.Range("E9").Value = 0.0388214940135581
$ws.Range("F9").Value = 0.07644845162506543
$ws.Range("G9").Value = 0.1442422934923629
$ws.Range("H9").Value = 0.2721551423170391
$ws.Range("I9").Value = 0.5523215202283346
$ws.Range("J9").Value = 1.045315968430635
$ws.Range("K9").Value = 1.972293349381438
$ws.Range("L9").Value = 4.460708233363363
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0.4723196494102408
$ws.Range("O9").Value = 2.917726550808414
$ws.Range("P9").Value = 0.02441902216341692
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.04976275060681323
$ws.Range("I10").Value = 0.09389193797172701
$ws.Range("J10").Value = 0.1771545163518279
$ws.Range("K10").Value = 0.3342536467113974
$ws.Range("L10").Value = 0.4887197018039945
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0.2284444358209828
$ws.Range("O10").Value = 0.03842101081093909
$ws.Range("P10").Value = 0.08189411429438404
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.1559232852346815
$ws.Range("J11").Value = 0.2941947389780779
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0.0806804875513108
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0.02892281900020256
$ws.Range("P11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("E13").Value = 0.1205364403587254
$ws.Range("F13").Value = 0.2274271386426275
$ws.Range("G13").Value = 0.4291076062744276
$ws.Range("H13").Value = 0.8096366109231623
$ws.Range("I13").Value = 1.648151966539643
$ws.Range("J13").Value = 3.109719224184783
$ws.Range("K13").Value = 5.438284387661469
$ws.Range("L13").Value = 1.600412105269389
$ws.Range("M13").Value = 0.6688125948725911
$ws.Range("N13").Value = 0.997955817856817
$ws.Range("O13").Value = 0.5059256759218019
$ws.Range("P13").Value = 0.5975130118787235
$ws.Range("E14").Value = 0.4103833652110969
$ws.Range("F14").Value = 0.815106777827997
$ws.Range("G14").Value = 1.537936590942427
$ws.Range("H14").Value = 2.951528534057
$ws.Range("I14").Value = 6.135225784982522
$ws.Range("J14").Value = 11.61669115209481
$ws.Range("K14").Value = 20.93408309243074
$ws.Range("L14").Value = 36.63585134629233
$ws.Range("M14").Value = 9.73681889192779
$ws.Range("N14").Value = 12.02981464179621
$ws.Range("O14").Value = 17.29672823998208
$ws.Range("P14").Value = 6.637279107260948

$ws = $wb.Worksheets.Item("INVESTMENT")
$ws.Range("E2").Value = 39.86214284453398
$ws.Range("F2").Value = 68.29554974478945
$ws.Range("G2").Value = 114.1792745700474
$ws.Range("H2").Value = 196.9668497107533
$ws.Range("I2").Value = 369.6343227191189
$ws.Range("J2").Value = 661.9607388769967
$ws.Range("K2").Value = 1278.98379610933
$ws.Range("L2").Value = 2797.023729318964
$ws.Range("M2").Value = 100.6217618475533
$ws.Range("N2").Value = 150.1408218215421
$ws.Range("O2").Value = 1632.326233186462
$ws.Range("P2").Value = 88.13221187467148
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("E5").Value = 110.9256046750557
$ws.Range("F5").Value = 190.0481161186478
$ws.Range("G5").Value = 317.7301612318728
$ws.Range("H5").Value = 548.105679875714
$ws.Range("I5").Value = 1028.592740640767
$ws.Range("J5").Value = 1880.666823139359
$ws.Range("K5").Value = 3526.35272568641
$ws.Range("L5").Value = 7799.312831696379
$ws.Range("M5").Value = 238.0170346810957
$ws.Range("N5").Value = 355.1525290193985
$ws.Range("O5").Value = 4625.08991953196
$ws.Range("P5").Value = 208.4734688116119
$ws.Range("E6").Value = 179.3796428004029
$ws.Range("F6").Value = 307.3299738515526
$ws.Range("G6").Value = 513.8067355652132
$ws.Range("H6").Value = 886.35082369839
$ws.Range("I6").Value = 1672.719661950729
$ws.Range("J6").Value = 3093.614056095011
$ws.Range("K6").Value = 5735.779357100485
$ws.Range("L6").Value = 12639.98476960104
$ws.Range("M6").Value = 19237.56615487933
$ws.Range("N6").Value = 12422.85913917319
$ws.Range("O6").Value = 2350.467959323159
$ws.Range("P6").Value = 11279.06340484421
$ws.Range("F7").Value = 138.2049736570903
$ws.Range("G7").Value = 231.0566895369785
$ws.Range("H7").Value = 423.2798143982842
$ws.Range("I7").Value = 693.2986735065595
$ws.Range("J7").Value = 1339.564097813024
$ws.Range("K7").Value = 2437.211157282198
$ws.Range("L7").Value = 5213.754755200701
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 2463.80722931935
$ws.Range("O7").Value = 88.44363220694997
$ws.Range("P7").Value = 865.9217718033158
$ws.Range("E8").Value = 688.0690180002125
$ws.Range("F8").Value = 1178.864167688006
$ws.Range("G8").Value = 1970.873006897621
$ws.Range("H8").Value = 3480.18283528256
$ws.Range("I8").Value = 6649.937750715355
$ws.Range("J8").Value = 12033.27555999994
$ws.Range("K8").Value = 22076.81428615932
$ws.Range("L8").Value = 46944.5066127752
$ws.Range("M8").Value = 656.3711408079164
$ws.Range("N8").Value = 9284.75129883386
$ws.Range("O8").Value = 25975.78195723822
$ws.Range("P8").Value = 729.4826208725113
$ws.Range("E9").Value = 176.8211217892161
$ws.Range("F9").Value = 316.1829368252248
$ws.Range("G9").Value = 528.6074787162196
$ws.Range("H9").Value = 911.8830909404908
$ws.Range("I9").Value = 1811.714552125527
$ws.Range("J9").Value = 3240.905900442397
$ws.Range("K9").Value = 5915.414481050653
$ws.Range("L9").Value = 12638.04348186237
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 1332.474426182715
$ws.Range("O9").Value = 8110.73828119956
$ws.Range("P9").Value = 63.92061096839087
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 132.7455529588759
$ws.Range("I10").Value = 230.8957526195319
$ws.Range("J10").Value = 389.3791860268568
$ws.Range("K10").Value = 708.4388855654279
$ws.Range("L10").Value = 991.5684790341954
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 457.2810852911591
$ws.Range("O10").Value = 80.07084652840788
$ws.Range("P10").Value = 160.7145408437063
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 347.6792524852219
$ws.Range("J11").Value = 622.6424352973517
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 155.5086234622787
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 54.65456110231143
$ws.Range("P11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("E13").Value = 533.7417673137425
$ws.Range("F13").Value = 914.4562940987495
$ws.Range("G13").Value = 1623.532473143449
$ws.Range("H13").Value = 2793.567253699799
$ws.Range("I13").Value = 5117.232912502095
$ws.Range("J13").Value = 9028.771118248318
$ws.Range("K13").Value = 15013.67173823803
$ws.Range("L13").Value = 4274.372535385646
$ws.Range("M13").Value = 1736.073843799446
$ws.Range("N13").Value = 2590.449112249075
$ws.Range("O13").Value = 1367.267377102996
$ws.Range("P13").Value = 1520.585855608614
$ws.Range("E14").Value = 1728.799297423164
$ws.Range("F14").Value = 3113.38201198406
$ws.Range("G14").Value = 5299.785819661402
$ws.Range("H14").Value = 9373.081900564866
$ws.Range("I14").Value = 17921.7056192649
$ws.Range("J14").Value = 32290.77991593925
$ws.Range("K14").Value = 56692.66642719185
$ws.Range("L14").Value = 93298.56719487447
$ws.Range("M14").Value = 22124.15855947763
$ws.Range("N14").Value = 29056.91564189029
$ws.Range("O14").Value = 44284.84076742003
$ws.Range("P14").Value = 14916.29448562704

$ws = $wb.Worksheets.Item("REMOVAL")
$ws.Range("B2").Value = 0.05031567005800005
$ws.Range("C2").Value = 0.1452508517849684
$ws.Range("D2").Value = 0.3243737516518166
$ws.Range("E2").Value = 0.6623413277711184
$ws.Range("F2").Value = 1.300015698791594
$ws.Range("G2").Value = 2.503174321655522
$ws.Range("H2").Value = 4.773283859624301
$ws.Range("I2").Value = 15.51391115206189
$ws.Range("J2").Value = 14.30514719062254
$ws.Range("K2").Value = 12.35679933717087
$ws.Range("L2").Value = 7.398472213920517
$ws.Range("M2").Value = 4.429738608402963
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("B5").Value = 0.1308207421508001
$ws.Range("C5").Value = 0.377652214640918
$ws.Range("D5").Value = 0.8433717542947233
$ws.Range("E5").Value = 1.722087452204907
$ws.Range("F5").Value = 3.380040816858143
$ws.Range("G5").Value = 6.508253236304355
$ws.Range("H5").Value = 12.41053803502318
$ws.Range("I5").Value = 40.33616899536089
$ws.Range("J5").Value = 37.66464050266554
$ws.Range("K5").Value = 32.24089844397629
$ws.Range("L5").Value = 19.30381685264177
$ws.Range("M5").Value = 18.7304032791313
$ws.Range("B6").Value = 0.2264205152610002
$ws.Range("C6").Value = 0.6536288330323581
$ws.Range("D6").Value = 1.459681882433175
$ws.Range("E6").Value = 2.980535974970033
$ws.Range("F6").Value = 5.850070644562171
$ws.Range("G6").Value = 11.26428444744985
$ws.Range("H6").Value = 21.47977736830935
$ws.Range("I6").Value = 69.81260018427847
$ws.Range("J6").Value = 143.8767053993182
$ws.Range("K6").Value = 184.9255694299553
$ws.Range("L6").Value = 136.1092443018842
$ws.Range("M6").Value = 106.831808521619
$ws.Range("C7").Value = 0.1710732781972001
$ws.Range("D7").Value = 0.4938528960688927
$ws.Range("E7").Value = 1.102870755616177
$ws.Range("F7").Value = 2.251960514421802
$ws.Range("G7").Value = 4.420053375891418
$ws.Range("H7").Value = 8.51079269362877
$ws.Range("I7").Value = 27.87378467668022
$ws.Range("J7").Value = 16.68906452230515
$ws.Range("K7").Value = 9.992359410836801
$ws.Range("L7").Value = 5.982794689414238
$ws.Range("M7").Value = 3.582120180431509
$ws.Range("B8").Value = 0.7346087828468006
$ws.Range("C8").Value = 2.12066243606054
$ws.Range("D8").Value = 4.735856774116524
$ws.Range("E8").Value = 9.67018338545833
$ws.Range("F8").Value = 18.98022920235727
$ws.Range("G8").Value = 36.5463450961706
$ws.Range("H8").Value = 69.68994435051479
$ws.Range("I8").Value = 226.5031028201035
$ws.Range("J8").Value = 219.1520205754426
$ws.Range("K8").Value = 206.1444281485774
$ws.Range("L8").Value = 125.2426837098748
$ws.Range("M8").Value = 121.1749002825676
$ws.Range("B9").Value = 0.1766377977616894
$ws.Range("C9").Value = 0.524478252655737
$ws.Range("D9").Value = 1.180780688045989
$ws.Range("E9").Value = 2.419086585588516
$ws.Range("F9").Value = 4.755511704865749
$ws.Range("G9").Value = 9.16385890633109
$ws.Range("H9").Value = 17.48149121062638
$ws.Range("I9").Value = 56.83563023669046
$ws.Range("J9").Value = 55.16367563797809
$ws.Range("K9").Value = 48.75765547434567
$ws.Range("L9").Value = 30.84942042198976
$ws.Range("M9").Value = 30.70049707438785
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.2264205152610002
$ws.Range("F10").Value = 0.6536288330323581
$ws.Range("G10").Value = 1.459681882433175
$ws.Range("H10").Value = 2.980535974970033
$ws.Range("I10").Value = 7.427885261386383
$ws.Range("J10").Value = 4.447349286416349
$ws.Range("K10").Value = 2.662792299472913
$ws.Range("L10").Value = 1.594312111213937
$ws.Range("M10").Value = 0.9545735536589105
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0.7094509478178006
$ws.Range("G11").Value = 2.048037010168055
$ws.Range("H11").Value = 2.048037010168055
$ws.Range("I11").Value = 1.226235410914942
$ws.Range("J11").Value = 0.7341924367169282
$ws.Range("K11").Value = 0.4395881323718608
$ws.Range("L11").Value = 0.2631976529018433
$ws.Range("M11").Value = 0.1575861571131749
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("B13").Value = 0.5484408036322005
$ws.Range("C13").Value = 1.583234284456156
$ws.Range("D13").Value = 3.535673893004802
$ws.Range("E13").Value = 7.21952047270519
$ws.Range("F13").Value = 14.17017111682837
$ws.Range("G13").Value = 27.28460010604518
$ws.Range("H13").Value = 52.02879406990486
$ws.Range("I13").Value = 63.75978570685265
$ws.Range("J13").Value = 38.17533894061589
$ws.Range("K13").Value = 22.85698559169206
$ws.Range("L13").Value = 13.68532159338543
$ws.Range("M13").Value = 8.193907563316486
$ws.Range("B14").Value = 1.867244311710491
$ws.Range("C14").Value = 5.575980150827878
$ws.Range("D14").Value = 12.57359163961592
$ws.Range("E14").Value = 26.00304646957527
$ws.Range("F14").Value = 52.05107947953525
$ws.Range("G14").Value = 101.1982883824492
$ws.Range("H14").Value = 191.4031945727697
$ws.Range("I14").Value = 509.2891044443294
$ws.Range("J14").Value = 530.2081344920812
$ws.Range("K14").Value = 520.377076268399
$ws.Range("L14").Value = 340.4292635472265
$ws.Range("M14").Value = 294.7555352206288
